$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.027.50"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.642.96"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.41"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5158"
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2600"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06394"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.90"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07781"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.318"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "1.655.82"
$ws.Range("E13").Value = "  -6.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5501"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.89"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "0.0₅7782"
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("D17").Value = "26.080.93"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "199.97"
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.492"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.04"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.137"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.905"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.39"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1220"
$ws.Range("E26").Value = "  +6.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.922"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04914"
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.320"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.252"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.550"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.384"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9207"
$ws.Range("E35").Value = "  +2.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.603"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5605"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("D38").Value = "1.114.91"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01575"
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.003"
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.545"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.593"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8128"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.92"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "0.0₈120"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").Value = "1.783.61"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4545"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.53"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05260"
$ws.Range("E50").Value = "  +3.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09615"
$ws.Range("E51").Value = "  +0.52%  "
